# Reviewed DigitalShares contract: the per-method gas table collapses
# initialize/addShare/distributeDividends/registerStock/unregisterStock
# into a single "distribute" row, and the Method/Snapshots/gas comparison
# table keeps only "withdraw" (with refreshed gas numbers) and drops
# "sendShares" entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Delete the rows that disappear, bottom-up (by their ORIGINAL row
# numbers) so every subsequent Delete() still targets the right row. ---
$ws.Rows.Item(10).Delete()   # sendShares snapshot row
$ws.Rows.Item(7).Delete()    # unregisterStock
$ws.Rows.Item(6).Delete()    # registerStock
$ws.Rows.Item(5).Delete()    # distributeDividends
$ws.Rows.Item(4).Delete()    # addShare

# After the deletes the surviving rows are (old -> new):
#   1 -> 1 (DigitalShare)
#   2 -> 2 (Method/gas header, style preserved)
#   3 -> 3 (was "initialize")
#   9 -> 5 (Method/Snapshots/gas header row, style preserved)
#   11 -> 6 (was "withdraw")

# --- Row 3: initialize -> distribute -------------------------------------
$ws.Range("A3").Value = "distribute"
$ws.Range("B3").Value = 74320

# --- Row 6 (was row 11): refresh the withdraw gas numbers and drop the
# trailing F column that no longer exists. ---------------------------------
$ws.Range("B6").Value = 69511
$ws.Range("C6").Value = 80482
$ws.Range("D6").Value = 190192
$ws.Range("E6").Value = 1287292
$ws.Range("F6").Clear()

# --- Column E width / selection tidy-up -----------------------------------
$ws.Columns.Item(5).ColumnWidth = 10.666666666666666
$ws.Range("E6").Select()
